# Add a new "23-jun" data column (L), matching the look/format of the
# existing "22-jun" column (K), and remove the picture that used to be
# anchored over the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the embedded picture anchored over the sheet.
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Clone the formatting of column K (header + data) into column L so the
# new column keeps the same number format / alignment / style.
$ws.Range("K1:K11").Copy() | Out-Null
$ws.Range("L1:L11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header for the new column.
$ws.Range("L1").Value = "23-jun"

# New data values for column L ("23-jun").
$ws.Range("L2").Value = 13
$ws.Range("L3").Value = 14
$ws.Range("L4").Value = 8
$ws.Range("L5").Value = 10
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 16
$ws.Range("L8").Value = 14
$ws.Range("L9").Value = 13
$ws.Range("L10").Value = 13
$ws.Range("L11").Value = 6

# Update the selection to mirror where the user ended up after the edit.
[void]$ws.Range("M7").Select()
